# Update "Estado de Cuenta" workbook:
#  - refresh VALOR MORA total and Cant. Trabajadores
#  - replace the worker detail table (remove DAGOBERTO BOLAÑOS ORTIZ,
#    regroup rows by period instead of by worker, and refresh the
#    Valor Mora / Salario Basico figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- header totals -------------------------------------------------
$ws.Range("E11").Value = 573256
$ws.Range("C13").Value = 4

# ---- worker detail table (rows 16-31) -------------------------------
# columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo, F=Valor Mora, G=Salario Basico
$rows = @(
    @("CC","73574401","HEMITT ENRIQUE ROCHA CASTRO","1911",33125,908526),
    @("CC","73581603","REINALDO AVILA RICARDO","1911",46805,908526),
    @("CC","1143372467","YERIS ANDRES NAVARRO MERCADO","1911",33125,1000000),
    @("CC","73574315","SANDER ORTEGA MARRUGO","1911",33125,1000000),
    @("CC","73574401","HEMITT ENRIQUE ROCHA CASTRO","1912",33125,908526),
    @("CC","73581603","REINALDO AVILA RICARDO","1912",46805,908526),
    @("CC","1143372467","YERIS ANDRES NAVARRO MERCADO","1912",33125,1000000),
    @("CC","73574315","SANDER ORTEGA MARRUGO","1912",33125,1000000),
    @("CC","73574401","HEMITT ENRIQUE ROCHA CASTRO","2002",35112,908526),
    @("CC","73581603","REINALDO AVILA RICARDO","2002",35112,908526),
    @("CC","1143372467","YERIS ANDRES NAVARRO MERCADO","2002",35112,1000000),
    @("CC","73574315","SANDER ORTEGA MARRUGO","2002",35112,1000000),
    @("CC","73574401","HEMITT ENRIQUE ROCHA CASTRO","2003",35112,908526),
    @("CC","73581603","REINALDO AVILA RICARDO","2003",35112,908526),
    @("CC","1143372467","YERIS ANDRES NAVARRO MERCADO","2003",35112,1000000),
    @("CC","73574315","SANDER ORTEGA MARRUGO","2003",35112,1000000)
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r++
}

# rows 32-35 used to hold the 5th worker (now removed); delete them so the
# signature block below shifts up from rows 40/41 to rows 36/37.
$ws.Range("B32:B35").EntireRow.Delete()
